$wb = $excel.ActiveWorkbook

# --- Sheet1 (Neg_Change) ---
$ws1 = $wb.Worksheets.Item("Neg_Change")

$ws1.Cells.Item(2, 1).Value = "BAJAJ-AUTO"
$ws1.Cells.Item(2, 2).Value = 9520
$ws1.Cells.Item(2, 3).Value = 9684
$ws1.Cells.Item(2, 4).Value = 9465
$ws1.Cells.Item(2, 5).Value = 9499
$ws1.Cells.Item(2, 6).Value = 304245
$ws1.Cells.Item(2, 7).Value = 733651
$ws1.Cells.Item(2, 8).Value = -0.5853000950043004
$ws1.Cells.Item(2, 9).Value = "BAJAJ-AUTO"

$ws1.Cells.Item(3, 1).Value = "ITC"
$ws1.Cells.Item(3, 2).Value = 350.2
$ws1.Cells.Item(3, 3).Value = 353.95
$ws1.Cells.Item(3, 4).Value = 346.05
$ws1.Cells.Item(3, 5).Value = 349.6
$ws1.Cells.Item(3, 6).Value = 48358124
$ws1.Cells.Item(3, 7).Value = 118790358
$ws1.Cells.Item(3, 8).Value = -0.5929120442586763
$ws1.Cells.Item(3, 9).Value = "ITC"

$ws1.Cells.Item(4, 1).Value = "COALINDIA"
$ws1.Cells.Item(4, 2).Value = 430.95
$ws1.Cells.Item(4, 3).Value = 436.7
$ws1.Cells.Item(4, 4).Value = 423.5
$ws1.Cells.Item(4, 5).Value = 426.85
$ws1.Cells.Item(4, 6).Value = 17258568
$ws1.Cells.Item(4, 7).Value = 35090129
$ws1.Cells.Item(4, 8).Value = -0.5081645895345668
$ws1.Cells.Item(4, 9).Value = "COALINDIA"

$ws1.Cells.Item(5, 1).Value = "JINDALSTEL"
$ws1.Cells.Item(5, 2).Value = 1084
$ws1.Cells.Item(5, 3).Value = 1087.5
$ws1.Cells.Item(5, 4).Value = 1074.4
$ws1.Cells.Item(5, 5).Value = 1081
$ws1.Cells.Item(5, 6).Value = 900540
$ws1.Cells.Item(5, 7).Value = 1931379
$ws1.Cells.Item(5, 8).Value = -0.5337321157577047
$ws1.Cells.Item(5, 9).Value = "JINDALSTEL"

$ws1.Cells.Item(6, 1).Value = "UPL"
$ws1.Cells.Item(6, 2).Value = 810
$ws1.Cells.Item(6, 3).Value = 811.45
$ws1.Cells.Item(6, 4).Value = 798.1
$ws1.Cells.Item(6, 5).Value = 802.8
$ws1.Cells.Item(6, 6).Value = 2068765
$ws1.Cells.Item(6, 7).Value = 4562526
$ws1.Cells.Item(6, 8).Value = -0.5465746386979493
$ws1.Cells.Item(6, 9).Value = "UPL"

$ws1.Cells.Item(7, 1).Value = "INDUSINDBK"
$ws1.Cells.Item(7, 2).Value = 905
$ws1.Cells.Item(7, 3).Value = 910
$ws1.Cells.Item(7, 4).Value = 892.5
$ws1.Cells.Item(7, 5).Value = 897.85
$ws1.Cells.Item(7, 6).Value = 2521435
$ws1.Cells.Item(7, 7).Value = 6250730
$ws1.Cells.Item(7, 8).Value = -0.5966175150742393
$ws1.Cells.Item(7, 9).Value = "INDUSINDBK"

$ws1.Cells.Item(8, 1).Value = "PETRONET"
$ws1.Cells.Item(8, 2).Value = 291.05
$ws1.Cells.Item(8, 3).Value = 291.6
$ws1.Cells.Item(8, 4).Value = 286.55
$ws1.Cells.Item(8, 5).Value = 287.95
$ws1.Cells.Item(8, 6).Value = 1247137
$ws1.Cells.Item(8, 7).Value = 2451474
$ws1.Cells.Item(8, 8).Value = -0.4912705580397753
$ws1.Cells.Item(8, 9).Value = "PETRONET"

$ws1.Cells.Item(9, 1).Value = "IREDA"
$ws1.Cells.Item(9, 2).Value = 147.5
$ws1.Cells.Item(9, 3).Value = 149.48
$ws1.Cells.Item(9, 4).Value = 143.64
$ws1.Cells.Item(9, 5).Value = 144.27
$ws1.Cells.Item(9, 6).Value = 12236304
$ws1.Cells.Item(9, 7).Value = 28577864
$ws1.Cells.Item(9, 8).Value = -0.5718258019563673
$ws1.Cells.Item(9, 9).Value = "IREDA"

$ws1.Cells.Item(10, 1).Value = "INDUSTOWER"
$ws1.Cells.Item(10, 2).Value = 439.3
$ws1.Cells.Item(10, 3).Value = 442.65
$ws1.Cells.Item(10, 4).Value = 432
$ws1.Cells.Item(10, 5).Value = 434
$ws1.Cells.Item(10, 6).Value = 5894990
$ws1.Cells.Item(10, 7).Value = 11982152
$ws1.Cells.Item(10, 8).Value = -0.5080190937320775
$ws1.Cells.Item(10, 9).Value = "INDUSTOWER"

$ws1.Cells.Item(11, 1).Value = "IIFL"
$ws1.Cells.Item(11, 2).Value = 645
$ws1.Cells.Item(11, 3).Value = 649.1
$ws1.Cells.Item(11, 4).Value = 638.15
$ws1.Cells.Item(11, 5).Value = 647.75
$ws1.Cells.Item(11, 6).Value = 1101439
$ws1.Cells.Item(11, 7).Value = 2441538
$ws1.Cells.Item(11, 8).Value = -0.5488749304741519
$ws1.Cells.Item(11, 9).Value = "IIFL"

$ws1.Cells.Item(12, 1).Value = "MCX"
$ws1.Cells.Item(12, 2).Value = 2225
$ws1.Cells.Item(12, 3).Value = 2232
$ws1.Cells.Item(12, 4).Value = 2176
$ws1.Cells.Item(12, 5).Value = 2198
$ws1.Cells.Item(12, 6).Value = 1638124
$ws1.Cells.Item(12, 7).Value = 3418580
$ws1.Cells.Item(12, 8).Value = -0.5208174154180976
$ws1.Cells.Item(12, 9).Value = "MCX"

$ws1.Cells.Item(13, 1).Value = "ABFRL"
$ws1.Cells.Item(13, 2).Value = 77.9
$ws1.Cells.Item(13, 3).Value = 77.93
$ws1.Cells.Item(13, 4).Value = 76.77
$ws1.Cells.Item(13, 5).Value = 77.05
$ws1.Cells.Item(13, 6).Value = 1677318
$ws1.Cells.Item(13, 7).Value = 3695684
$ws1.Cells.Item(13, 8).Value = -0.5461413908764927
$ws1.Cells.Item(13, 9).Value = "ABFRL"

$ws1.Cells.Item(14, 1).Value = "INOXWIND"
$ws1.Cells.Item(14, 2).Value = 129
$ws1.Cells.Item(14, 3).Value = 129
$ws1.Cells.Item(14, 4).Value = 125.27
$ws1.Cells.Item(14, 5).Value = 125.66
$ws1.Cells.Item(14, 6).Value = 4493882
$ws1.Cells.Item(14, 7).Value = 9949850
$ws1.Cells.Item(14, 8).Value = -0.5483467589963668
$ws1.Cells.Item(14, 9).Value = "INOXWIND"

# --- Sheet2 (Pos_Change) ---
$ws2 = $wb.Worksheets.Item("Pos_Change")

$ws2.Cells.Item(2, 1).Value = "ULTRACEMCO"
$ws2.Cells.Item(2, 2).Value = 11931
$ws2.Cells.Item(2, 3).Value = 12114
$ws2.Cells.Item(2, 4).Value = 11860
$ws2.Cells.Item(2, 5).Value = 12059
$ws2.Cells.Item(2, 6).Value = 302803
$ws2.Cells.Item(2, 7).Value = 206244
$ws2.Cells.Item(2, 8).Value = 0.4681784682221058
$ws2.Cells.Item(2, 9).Value = "ULTRACEMCO"

$ws2.Cells.Item(3, 1).Value = "TATACONSUM"
$ws2.Cells.Item(3, 2).Value = 1165
$ws2.Cells.Item(3, 3).Value = 1191
$ws2.Cells.Item(3, 4).Value = 1165
$ws2.Cells.Item(3, 5).Value = 1183
$ws2.Cells.Item(3, 6).Value = 1521762
$ws2.Cells.Item(3, 7).Value = 955406
$ws2.Cells.Item(3, 8).Value = 0.5927909182065006
$ws2.Cells.Item(3, 9).Value = "TATACONSUM"

$ws2.Cells.Item(4, 1).Value = "SBIN"
$ws2.Cells.Item(4, 2).Value = 1000
$ws2.Cells.Item(4, 3).Value = 1015.5
$ws2.Cells.Item(4, 4).Value = 1000
$ws2.Cells.Item(4, 5).Value = 1005.4
$ws2.Cells.Item(4, 6).Value = 11259060
$ws2.Cells.Item(4, 7).Value = 7357435
$ws2.Cells.Item(4, 8).Value = 0.5302969037443076
$ws2.Cells.Item(4, 9).Value = "SBIN"

$ws2.Cells.Item(5, 1).Value = "CIPLA"
$ws2.Cells.Item(5, 2).Value = 1511
$ws2.Cells.Item(5, 3).Value = 1525
$ws2.Cells.Item(5, 4).Value = 1507.2
$ws2.Cells.Item(5, 5).Value = 1520.1
$ws2.Cells.Item(5, 6).Value = 1532549
$ws2.Cells.Item(5, 7).Value = 964611
$ws2.Cells.Item(5, 8).Value = 0.5887741276016965
$ws2.Cells.Item(5, 9).Value = "CIPLA"

$ws2.Cells.Item(6, 1).Value = "SUNPHARMA"
$ws2.Cells.Item(6, 2).Value = 1728
$ws2.Cells.Item(6, 3).Value = 1746
$ws2.Cells.Item(6, 4).Value = 1723.4
$ws2.Cells.Item(6, 5).Value = 1725.7
$ws2.Cells.Item(6, 6).Value = 1416698
$ws2.Cells.Item(6, 7).Value = 974123
$ws2.Cells.Item(6, 8).Value = 0.4543317425006904
$ws2.Cells.Item(6, 9).Value = "SUNPHARMA"

$ws2.Cells.Item(7, 1).Value = "ETERNAL"
$ws2.Cells.Item(7, 2).Value = 280.5
$ws2.Cells.Item(7, 3).Value = 284
$ws2.Cells.Item(7, 4).Value = 279.6
$ws2.Cells.Item(7, 5).Value = 281.9
$ws2.Cells.Item(7, 6).Value = 19328131
$ws2.Cells.Item(7, 7).Value = 13624798
$ws2.Cells.Item(7, 8).Value = 0.4185994537313507
$ws2.Cells.Item(7, 9).Value = "ETERNAL"

$ws2.Cells.Item(8, 1).Value = "TECHM"
$ws2.Cells.Item(8, 2).Value = 1601.9
$ws2.Cells.Item(8, 3).Value = 1610.1
$ws2.Cells.Item(8, 4).Value = 1576.1
$ws2.Cells.Item(8, 5).Value = 1598.8
$ws2.Cells.Item(8, 6).Value = 907810
$ws2.Cells.Item(8, 7).Value = 637080
$ws2.Cells.Item(8, 8).Value = 0.4249544798141521
$ws2.Cells.Item(8, 9).Value = "TECHM"

$ws2.Cells.Item(9, 1).Value = "BANKBARODA"
$ws2.Cells.Item(9, 2).Value = 308.8
$ws2.Cells.Item(9, 3).Value = 311.8
$ws2.Cells.Item(9, 4).Value = 304.65
$ws2.Cells.Item(9, 5).Value = 307
$ws2.Cells.Item(9, 6).Value = 16328380
$ws2.Cells.Item(9, 7).Value = 10317519
$ws2.Cells.Item(9, 8).Value = 0.5825878294966067
$ws2.Cells.Item(9, 9).Value = "BANKBARODA"

$ws2.Cells.Item(10, 1).Value = "PNB"
$ws2.Cells.Item(10, 2).Value = 126
$ws2.Cells.Item(10, 3).Value = 128.24
$ws2.Cells.Item(10, 4).Value = 124.38
$ws2.Cells.Item(10, 5).Value = 125
$ws2.Cells.Item(10, 6).Value = 19053073
$ws2.Cells.Item(10, 7).Value = 13205239
$ws2.Cells.Item(10, 8).Value = 0.4428419659803204
$ws2.Cells.Item(10, 9).Value = "PNB"

$ws2.Cells.Item(11, 1).Value = "ABB"
$ws2.Cells.Item(11, 2).Value = 5230
$ws2.Cells.Item(11, 3).Value = 5249.5
$ws2.Cells.Item(11, 4).Value = 5131
$ws2.Cells.Item(11, 5).Value = 5165.5
$ws2.Cells.Item(11, 6).Value = 101498
$ws2.Cells.Item(11, 7).Value = 63490
$ws2.Cells.Item(11, 8).Value = 0.5986454559773192
$ws2.Cells.Item(11, 9).Value = "ABB"

$ws2.Cells.Item(12, 1).Value = "BPCL"
$ws2.Cells.Item(12, 2).Value = 380
$ws2.Cells.Item(12, 3).Value = 385.45
$ws2.Cells.Item(12, 4).Value = 374.4
$ws2.Cells.Item(12, 5).Value = 377.6
$ws2.Cells.Item(12, 6).Value = 5243413
$ws2.Cells.Item(12, 7).Value = 3614573
$ws2.Cells.Item(12, 8).Value = 0.4506313747156303
$ws2.Cells.Item(12, 9).Value = "BPCL"

$ws2.Cells.Item(13, 1).Value = "VOLTAS"
$ws2.Cells.Item(13, 2).Value = 1429.5
$ws2.Cells.Item(13, 3).Value = 1494.6
$ws2.Cells.Item(13, 4).Value = 1423.2
$ws2.Cells.Item(13, 5).Value = 1476
$ws2.Cells.Item(13, 6).Value = 1588073
$ws2.Cells.Item(13, 7).Value = 1109054
$ws2.Cells.Item(13, 8).Value = 0.4319167506721945
$ws2.Cells.Item(13, 9).Value = "VOLTAS"

$ws2.Cells.Item(14, 1).Value = "KALYANKJIL"
$ws2.Cells.Item(14, 2).Value = 492.7
$ws2.Cells.Item(14, 3).Value = 505.9
$ws2.Cells.Item(14, 4).Value = 488.65
$ws2.Cells.Item(14, 5).Value = 504.5
$ws2.Cells.Item(14, 6).Value = 2659371
$ws2.Cells.Item(14, 7).Value = 1727423
$ws2.Cells.Item(14, 8).Value = 0.53950190543949
$ws2.Cells.Item(14, 9).Value = "KALYANKJIL"

$ws2.Cells.Item(15, 1).Value = "GMRAIRPORT"
$ws2.Cells.Item(15, 2).Value = 106
$ws2.Cells.Item(15, 3).Value = 108.02
$ws2.Cells.Item(15, 4).Value = 104.7
$ws2.Cells.Item(15, 5).Value = 105.52
$ws2.Cells.Item(15, 6).Value = 11284687
$ws2.Cells.Item(15, 7).Value = 7760926
$ws2.Cells.Item(15, 8).Value = 0.4540387319760554
$ws2.Cells.Item(15, 9).Value = "GMRAIRPORT"

$ws2.Cells.Item(16, 1).Value = "IGL"
$ws2.Cells.Item(16, 2).Value = 194.33
$ws2.Cells.Item(16, 3).Value = 194.8
$ws2.Cells.Item(16, 4).Value = 190.31
$ws2.Cells.Item(16, 5).Value = 190.7
$ws2.Cells.Item(16, 6).Value = 1163937
$ws2.Cells.Item(16, 7).Value = 727464
$ws2.Cells.Item(16, 8).Value = 0.5999925769522616
$ws2.Cells.Item(16, 9).Value = "IGL"

$ws2.Cells.Item(17, 1).Value = "BANDHANBNK"
$ws2.Cells.Item(17, 2).Value = 147
$ws2.Cells.Item(17, 3).Value = 150.19
$ws2.Cells.Item(17, 4).Value = 145.95
$ws2.Cells.Item(17, 5).Value = 146.85
$ws2.Cells.Item(17, 6).Value = 12619074
$ws2.Cells.Item(17, 7).Value = 8134028
$ws2.Cells.Item(17, 8).Value = 0.551392987582536
$ws2.Cells.Item(17, 9).Value = "BANDHANBNK"

$ws2.Cells.Item(18, 1).Value = "NUVAMA"
$ws2.Cells.Item(18, 2).Value = 1475
$ws2.Cells.Item(18, 3).Value = 1505.2
$ws2.Cells.Item(18, 4).Value = 1467.5
$ws2.Cells.Item(18, 5).Value = 1489.8
$ws2.Cells.Item(18, 6).Value = 357422
$ws2.Cells.Item(18, 7).Value = 251732
$ws2.Cells.Item(18, 8).Value = 0.419851270398678
$ws2.Cells.Item(18, 9).Value = "NUVAMA"

$ws2.Cells.Item(19, 1).Value = "CDSL"
$ws2.Cells.Item(19, 2).Value = 1465
$ws2.Cells.Item(19, 3).Value = 1479
$ws2.Cells.Item(19, 4).Value = 1452.1
$ws2.Cells.Item(19, 5).Value = 1467.8
$ws2.Cells.Item(19, 6).Value = 1534959
$ws2.Cells.Item(19, 7).Value = 1061313
$ws2.Cells.Item(19, 8).Value = 0.4462830475081338
$ws2.Cells.Item(19, 9).Value = "CDSL"

$ws2.Cells.Item(20, 1).Value = "TATACHEM"
$ws2.Cells.Item(20, 2).Value = 755.25
$ws2.Cells.Item(20, 3).Value = 755.4
$ws2.Cells.Item(20, 4).Value = 745
$ws2.Cells.Item(20, 5).Value = 746.8
$ws2.Cells.Item(20, 6).Value = 442421
$ws2.Cells.Item(20, 7).Value = 298532
$ws2.Cells.Item(20, 8).Value = 0.4819885305427894
$ws2.Cells.Item(20, 9).Value = "TATACHEM"

$ws2.Cells.Item(21, 1).Value = "RBLBANK"
$ws2.Cells.Item(21, 2).Value = 320.9
$ws2.Cells.Item(21, 3).Value = 320.9
$ws2.Cells.Item(21, 4).Value = 309
$ws2.Cells.Item(21, 5).Value = 316.5
$ws2.Cells.Item(21, 6).Value = 8916084
$ws2.Cells.Item(21, 7).Value = 6278560
$ws2.Cells.Item(21, 8).Value = 0.420084223133967
$ws2.Cells.Item(21, 9).Value = "RBLBANK"

$ws2.Cells.Item(22, 1).Value = "NBCC"
$ws2.Cells.Item(22, 2).Value = 122.7
$ws2.Cells.Item(22, 3).Value = 124.27
$ws2.Cells.Item(22, 4).Value = 118.52
$ws2.Cells.Item(22, 5).Value = 119.65
$ws2.Cells.Item(22, 6).Value = 12521716
$ws2.Cells.Item(22, 7).Value = 8357400
$ws2.Cells.Item(22, 8).Value = 0.498278890564051
$ws2.Cells.Item(22, 9).Value = "NBCC"

# Remove now-unused trailing rows in Pos_Change (was 25 rows of data, now 22)
$ws2.Range("A23:I25").ClearContents()
